# 20150723 - Code cleanup, fomatting cleanup, add SPI1 and SPI2 pin definitions.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pins")

# --- Add SPI1 / SPI2 / SD Card pin definitions -----------------------------
# (New shared strings are created in the same order the original author's
#  workbook ended up with: "SD Card CS" (407), "SP12 CS" (408), "SPI2 CS" (409))

# PE10 (row 12): add SD Card CS in columns J and K
$ws.Range("J12").Value = "SD Card CS"
$ws.Range("K12").Value = "SD Card CS"

# PA3 (row 5): add SPI1-named column B entry ("SP12 CS") and column C entry ("SPI2 CS")
$ws.Range("B5").Value = "SP12 CS"
$ws.Range("C5").Value = "SPI2 CS"

# PA5 (row 7) SPI1 SCK -> mirror value into column C
$ws.Range("C7").Value = $ws.Range("B7").Value2

# PA6 (row 8) SPI1 MISO -> mirror value into column C
$ws.Range("C8").Value = $ws.Range("B8").Value2

# PA7 (row 9) SPI1 MOSI -> mirror value into column C
$ws.Range("C9").Value = $ws.Range("B9").Value2

# PB13 (row 32) SPI2 SCK -> mirror value into column C
$ws.Range("C32").Value = $ws.Range("B32").Value2

# PB14 (row 33) SPI2 MISO -> mirror value into column C
$ws.Range("C33").Value = $ws.Range("B33").Value2

# PB15 (row 34) SPI2 MOSI -> mirror value into column C
$ws.Range("C34").Value = $ws.Range("B34").Value2

# --- Update view / selection state -----------------------------------------
$ws.Activate()
$ws.Range("G26").Select()
